$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.101.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.08%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.294.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'316.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.35%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'104.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.44%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.50%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'39.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.63%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.31%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'8.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.74%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.963"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.89%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.31%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.643.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.88%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.295.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.119.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.91%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.33%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.08%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'PancakeSwap"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'3.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'Litecoin"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'73.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'277.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.39%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'10.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +9.96%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -3.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.73%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Toncoin"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +7.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Cosmos"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'10.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.92%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'22.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.67%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'36.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.48%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'163.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.16%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.83%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.84%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.76%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +3.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -6.76%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.76%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.94%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.98%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.34%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'100.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'69.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.23%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Algorand"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.225"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.76%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'FirstDigitalUSD"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'12.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.11%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'111.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.30%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'76.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.93%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -5.31%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.604.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.29%  "
$ws.Range("E51").Style = "Normal"
